# Apply the "employee-shift" StructureDefinition refresh (v5.0.0 -> v6.0.0)
# to the "Metadata" and "Elements" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Metadata": update version/date/publisher, replace the duplicated
# "Contact" row with a "Jurisdiction" row, and drop the extra duplicate
# row that used to follow it (21 rows -> 20 rows).
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact | No display for ContactDetail" row;
# remove it entirely so everything below shifts up by one row.
$meta.Rows.Item(11).Delete()

# ---------------------------------------------------------------------
# Sheet "Elements": the root Extension element's Short/Definition text
# now reflects the EmployeeShift extension itself.
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Employee Shift"
$elements.Range("L2").Value = "Code indicating the regular shift which the employee is scheduled to work (e.g., day, afternoon, midnight)"
